$d = $word.ActiveDocument

# Locate the "Proposal A:" heading paragraph (the one that currently reads
# "Proposal A:" + a manual line break + "Neurodiverse-friendly Planner application").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Proposal A:*Neurodiverse-friendly Planner application*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Proposal A' heading paragraph"
}

$r = $target.Range
$start = $r.Start
$end = $r.End

# Exclude the trailing paragraph mark from the replacement range.
$full = $d.Range($start, $end - 1)

# Rebuild the paragraph's run content so that "Proposal A" + ":" + the manual
# line break are removed, while the <w:lastRenderedPageBreak/> marker that sat
# on the very first run is preserved on the surviving run.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p>' +
       '<w:r><w:lastRenderedPageBreak/><w:t>Neurodiverse-friendly Planner application</w:t></w:r>' +
       '</w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xml)
